$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 525
$ws.Range("I38").Value = 31.666666
$ws.Range("J38").Value = 2005
$ws.Range("K38").Value = 94.99999800000001
$ws.Range("L38").Value = 6015
$ws.Range("M38").Value = 277.000002
$ws.Range("N38").Value = -6759

$ws.Range("H40").Value = 4746.25
$ws.Range("J40").Value = 4999
$ws.Range("L40").Value = 4999
$ws.Range("N40").Value = -5349

$ws.Range("H64").Value = 9460.352999999999
$ws.Range("I64").Value = 8566
$ws.Range("K64").Value = 8566
$ws.Range("M64").Value = -8318

$ws.Range("H67").Value = 9460.352999999999
$ws.Range("I67").Value = 8566
$ws.Range("K67").Value = 8566
$ws.Range("M67").Value = -7708

$ws.Range("H76").Value = 2000
$ws.Range("J76").Value = 2000
$ws.Range("L76").Value = 2000
$ws.Range("N76").Value = -2630

$ws.Range("H79").Value = 2000
$ws.Range("J79").Value = 2000
$ws.Range("L79").Value = 2000
$ws.Range("N79").Value = -4184

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("N86").Value = 0

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 0
$ws.Range("N89").Value = 0

$ws.Range("H95").Value = 31622.5
$ws.Range("J95").Value = 31622.5
$ws.Range("L95").Value = 31622.5
$ws.Range("N95").Value = -37114.5

$ws.Range("H96").Value = 1067.1
$ws.Range("I96").Value = 646.2
$ws.Range("J96").Value = 1488
$ws.Range("K96").Value = 1938.6
$ws.Range("L96").Value = 4464
$ws.Range("M96").Value = -565.6000000000001
$ws.Range("N96").Value = -7210

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 2200
$ws.Range("J21").Value = 2200
$ws.Range("L21").Value = 2200
$ws.Range("N21").Value = -2948

$ws.Range("H26").Value = 6500.8335
$ws.Range("I26").Value = 6701
$ws.Range("K26").Value = 6701
$ws.Range("M26").Value = -6371

$ws.Range("H27").Value = 7000
$ws.Range("J27").Value = 7000
$ws.Range("L27").Value = 7000
$ws.Range("N27").Value = -7368

$ws.Range("H28").Value = 2628.1667
$ws.Range("I28").Value = 2628.1667
$ws.Range("K28").Value = 2628.1667
$ws.Range("M28").Value = -2436.1667

$ws.Range("H30").Value = 10421.143
$ws.Range("I30").Value = 1999.5
$ws.Range("K30").Value = 1999.5
$ws.Range("M30").Value = -1849.5

$ws.Range("H31").Value = 5213.8
$ws.Range("I31").Value = 5213.8
$ws.Range("K31").Value = 5213.8
$ws.Range("M31").Value = -4919.8

$ws.Range("H76").Value = 28164.834
$ws.Range("J76").Value = 28164.834
$ws.Range("L76").Value = 28164.834
$ws.Range("N76").Value = -28840.834

$ws.Range("H79").Value = 28164.834
$ws.Range("J79").Value = 28164.834
$ws.Range("L79").Value = 28164.834
$ws.Range("N79").Value = -30504.834

$ws.Range("H99").Value = 2628.1667
$ws.Range("I99").Value = 2628.1667
$ws.Range("K99").Value = 2628.1667
$ws.Range("M99").Value = 366.8332999999998

$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 11714
$ws.Range("I97").Value = 11714
$ws.Range("K97").Value = 11714
$ws.Range("M97").Value = -10723

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4538
$ws.Range("I62").Value = 4538
$ws.Range("K62").Value = 4538
$ws.Range("M62").Value = -3914

$ws.Range("H65").Value = 4538
$ws.Range("I65").Value = 4538
$ws.Range("K65").Value = 22690
$ws.Range("M65").Value = -19570

$ws.Range("H69").Value = 9965.833000000001
$ws.Range("I69").Value = 9965.833000000001
$ws.Range("K69").Value = 9965.833000000001
$ws.Range("M69").Value = -9216.833000000001

$ws.Range("H72").Value = 9965.833000000001
$ws.Range("I72").Value = 9965.833000000001
$ws.Range("K72").Value = 29897.499
$ws.Range("M72").Value = -26153.499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3000
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6588

$ws.Range("H46").Value = 1918.9
$ws.Range("I46").Value = 449
$ws.Range("J46").Value = 2408.8667
$ws.Range("K46").Value = 1347
$ws.Range("L46").Value = 7226.6001
$ws.Range("M46").Value = -1256
$ws.Range("N46").Value = -7408.6001

$ws.Range("H88").Value = 2000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 2000
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4944.75
$ws.Range("I80").Value = 4599.6665
$ws.Range("J80").Value = 5980
$ws.Range("K80").Value = 4599.6665
$ws.Range("L80").Value = 5980
$ws.Range("M80").Value = -3601.6665
$ws.Range("N80").Value = -7976

$ws.Range("H83").Value = 4944.75
$ws.Range("I83").Value = 4599.6665
$ws.Range("J83").Value = 5980
$ws.Range("K83").Value = 22998.3325
$ws.Range("L83").Value = 29900
$ws.Range("M83").Value = -18006.3325
$ws.Range("N83").Value = -39884

$ws.Range("H99").Value = 29498
$ws.Range("I99").Value = 29498
$ws.Range("K99").Value = 29498
$ws.Range("M99").Value = -27252

$ws.Range("H126").Value = 1801.1428
$ws.Range("I126").Value = 1877.25
$ws.Range("J126").Value = 1699.6666
$ws.Range("K126").Value = 5631.75
$ws.Range("L126").Value = 5098.9998
$ws.Range("M126").Value = -3161.75
$ws.Range("N126").Value = -10038.9998

$ws.Range("H134").Value = 85000
$ws.Range("J134").Value = 85000
$ws.Range("L134").Value = 255000
$ws.Range("N134").Value = -260070

$ws.Range("H137").Value = 75000
$ws.Range("J137").Value = 75000
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 16666.666
$ws.Range("J26").Value = 16666.666
$ws.Range("L26").Value = 16666.666
$ws.Range("N26").Value = -17256.666

$ws.Range("H93").Value = 30305812
$ws.Range("I93").Value = 47621650
$ws.Range("J93").Value = 3099.5
$ws.Range("K93").Value = 47621650
$ws.Range("L93").Value = 3099.5
$ws.Range("M93").Value = -47620402
$ws.Range("N93").Value = -5595.5

$ws.Range("H132").Value = 2839.4
$ws.Range("I132").Value = 2924.25
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 8772.75
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -6242.75
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4925
$ws.Range("I62").Value = 4925
$ws.Range("K62").Value = 4925
$ws.Range("M62").Value = -4301

$ws.Range("H65").Value = 4925
$ws.Range("I65").Value = 4925
$ws.Range("K65").Value = 24625
$ws.Range("M65").Value = -21505
